# Daily attendance processing - 2025-11-05 13:57:03
# Normalize the "Recorded By" (column G) entries that still show the
# legacy ordering coming out of the sync job: "System" should be listed
# after the actual account(s) that recorded the session, not before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Find the last used row on the sheet so we cover every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
